$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.341.57"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.396.30"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "405.67"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.57"
$ws.Range("E6").Value = "  +15.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  +8.01%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.677"
$ws.Range("E9").Value = "  +9.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  +8.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.37"
$ws.Range("E11").Value = "  +9.41%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "3.935.03"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.54"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.73"
$ws.Range("E15").Value = "  +4.13%  "
$ws.Range("D16").Value = "3.380.09"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.51"
$ws.Range("E17").Value = "  +9.40%  "
$ws.Range("D18").Value = "61.162.83"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.03"
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000135"
$ws.Range("E20").Value = "  +17.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.26"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "83.03"
$ws.Range("E22").Value = "  +13.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.14"
$ws.Range("E23").Value = "  +6.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "308.05"
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.56"
$ws.Range("E26").Value = "  +15.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.73"
$ws.Range("E27").Value = "  +3.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.73"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.174"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.116"
$ws.Range("E31").Value = "  +6.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.80"
$ws.Range("E32").Value = "  +6.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.61"
$ws.Range("E33").Value = "  +6.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.64"
$ws.Range("E34").Value = "  +8.03%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0487"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.31"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("E39").Value = "  +5.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.03"
$ws.Range("E41").Value = "  +8.12%  "
$ws.Range("E42").Value = "  +5.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.16"
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.99"
$ws.Range("E44").Value = "  +6.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.286"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.02"
$ws.Range("E46").Value = "  +5.16%  "
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.77"
$ws.Range("E48").Value = "  +4.38%  "
$ws.Range("D49").Value = "2.146.15"
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("D50").Value = "3.716.43"
$ws.Range("E50").Value = "  +1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  +0.18%  "
